# Refresh the "ランサーズ" job-listing sheet with a newer scrape snapshot
# (2025-12-09 06:29:13 JST), replacing the old 23-row dataset with a new
# 13-row dataset, and tightening two column widths.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Drop all existing hyperlinks up front; we rebuild only the ones the
#     surviving rows need (F2:F13) further down. Doing this before the row
#     delete avoids leaving stale hyperlink entries pointing at refs that no
#     longer have a backing row. ---
$ws.Hyperlinks.Delete()

# --- Remove the now-stale rows 14-23 entirely (dimension collapses to
#     A1:H13 automatically). ---
$ws.Range("A14:H23").EntireRow.Delete()

# --- New snapshot timestamp applied to every data row. ---
$timestamp = "2025-12-09 06:29:13"
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 1).Value = $timestamp
}

# --- Row 2 ---
$ws.Cells.Item(2, 2).Value = "大手SIer等のAIソリューション開発・導入を支援してくださるエンジニア・PM募集"
$ws.Cells.Item(2, 4).Value = "200,000 円 ~ 300,000 円 / 固定"
$url2 = "https://www.lancers.jp/work/detail/5450158"
$ws.Cells.Item(2, 6).Value = $url2
$ws.Cells.Item(2, 7).Value = 368
$ws.Cells.Item(2, 8).Value = "🔥AI,Ai ◆開発"

# --- Row 3 ---
$ws.Cells.Item(3, 2).Value = "EC×AIプロダクト/業務改善リード"
$ws.Cells.Item(3, 4).Value = "200,000 円 ~ 300,000 円 / 固定"
$url3 = "https://www.lancers.jp/work/detail/5450024"
$ws.Cells.Item(3, 6).Value = $url3
$ws.Cells.Item(3, 7).Value = 338
$ws.Cells.Item(3, 8).Value = "🔥AI,Ai ◇業務改善"

# --- Row 4 ---
$ws.Cells.Item(4, 2).Value = "【せどり×ツール製作】APIを使用したせどりツールを製作できるエンジニアさんを募集します♪"
$ws.Cells.Item(4, 4).Value = "20,000 円 ~ 50,000 円 / 固定"
$url4 = "https://www.lancers.jp/work/detail/5217096"
$ws.Cells.Item(4, 6).Value = $url4
$ws.Cells.Item(4, 7).Value = 243
$ws.Cells.Item(4, 8).Value = "🔥API ◆ツール"

# --- Row 5 ---
$ws.Cells.Item(5, 2).Value = "【急募】pythonのコードのMac環境用インストーラー作成(Windows版は作成済)"
$ws.Cells.Item(5, 4).Value = "5,000 円 ~ 10,000 円 / 固定"
$url5 = "https://www.lancers.jp/work/detail/5442448"
$ws.Cells.Item(5, 6).Value = $url5
$ws.Cells.Item(5, 7).Value = 190
$ws.Cells.Item(5, 8).Value = "🔥Python"

# --- Row 6 ---
$ws.Cells.Item(6, 2).Value = "【バイナリ解析 / 逆コンパイル】EPCデータ解析ツール開発(継続依頼あり・高単価)"
$ws.Cells.Item(6, 4).Value = "200,000 円 ~ 300,000 円 / 固定"
$url6 = "https://www.lancers.jp/work/detail/5449973"
$ws.Cells.Item(6, 6).Value = $url6
$ws.Cells.Item(6, 7).Value = 128
$ws.Cells.Item(6, 8).Value = "◆ツール,開発"

# --- Row 7 ---
$ws.Cells.Item(7, 2).Value = "【単発/Stripeサブスクリプション実装】Laravel プラットフォーム開発エンジニア募集"
$ws.Cells.Item(7, 4).Value = "300,000 円 ~ 500,000 円 / 固定"
$url7 = "https://www.lancers.jp/work/detail/5449939"
$ws.Cells.Item(7, 6).Value = $url7
$ws.Cells.Item(7, 7).Value = 75
$ws.Cells.Item(7, 8).Value = "◆開発"

# --- Row 8 ---
$ws.Cells.Item(8, 2).Value = "Excel VBA一括自動処理ツール作成(データ転記・分類・置換・NGチェック)【エクセルマクロ】"
$ws.Cells.Item(8, 4).Value = "1,000 ~ 5,000 円 / 固定"
$url8 = "https://www.lancers.jp/work/detail/5450139"
$ws.Cells.Item(8, 6).Value = $url8
$ws.Cells.Item(8, 7).Value = 65
$ws.Cells.Item(8, 8).Value = "◆ツール"

# --- Row 9 ---
$ws.Cells.Item(9, 2).Value = "初回 【案件】Win2008(PHP5.3)→ Linux(LAMP)へのレガシー調査と移行"
$ws.Cells.Item(9, 4).Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$url9 = "https://www.lancers.jp/work/detail/5449999"
$ws.Cells.Item(9, 6).Value = $url9
$ws.Cells.Item(9, 7).Value = 40
$ws.Cells.Item(9, 8).Value = "○PHP"

# --- Row 10 ---
$ws.Cells.Item(10, 2).Value = "【急募】古いPHPとPerlプログラムのアップデート依頼"
$ws.Cells.Item(10, 4).Value = "100,000 円 ~ 200,000 円 / 固定"
$url10 = "https://www.lancers.jp/work/detail/5440861"
$ws.Cells.Item(10, 6).Value = $url10
$ws.Cells.Item(10, 7).Value = 33
$ws.Cells.Item(10, 8).Value = "○PHP"

# --- Row 11 (no skill-summary column in the new data) ---
$ws.Cells.Item(11, 2).Value = "【急募】社内システム保守運用・社内スタッフ教育まで依頼できる方を探しています!"
$ws.Cells.Item(11, 4).Value = "5,000 円 ~ 10,000 円 / 固定"
$url11 = "https://www.lancers.jp/work/detail/5449609"
$ws.Cells.Item(11, 6).Value = $url11
$ws.Cells.Item(11, 7).Value = 25
$ws.Cells.Item(11, 8).ClearContents()

# --- Row 12 (no skill-summary column in the new data) ---
$ws.Cells.Item(12, 2).Value = "X(旧ツイッター)自動ログインについて"
$ws.Cells.Item(12, 4).Value = "20,000 円 ~ 50,000 円 / 固定"
$url12 = "https://www.lancers.jp/work/detail/5449817"
$ws.Cells.Item(12, 6).Value = $url12
$ws.Cells.Item(12, 7).Value = 13
$ws.Cells.Item(12, 8).ClearContents()

# --- Row 13 (no skill-summary column in the new data) ---
$ws.Cells.Item(13, 2).Value = "【アカウント復活】削除したxのアカウントを再生させたい!"
$ws.Cells.Item(13, 4).Value = "5,000 円 ~ 10,000 円 / 固定"
$url13 = "https://www.lancers.jp/work/detail/5449948"
$ws.Cells.Item(13, 6).Value = $url13
$ws.Cells.Item(13, 7).Value = 10
$ws.Cells.Item(13, 8).ClearContents()

# --- Rebuild the F-column hyperlinks for the surviving rows, reusing the
#     exact URL strings assigned above (reading `.Value` back from a Range
#     in this host returns the property descriptor, not the stored value,
#     so we pass the known-good string instead of round-tripping it). ---
$ws.Hyperlinks.Add($ws.Cells.Item(2, 6), $url2)
$ws.Hyperlinks.Add($ws.Cells.Item(3, 6), $url3)
$ws.Hyperlinks.Add($ws.Cells.Item(4, 6), $url4)
$ws.Hyperlinks.Add($ws.Cells.Item(5, 6), $url5)
$ws.Hyperlinks.Add($ws.Cells.Item(6, 6), $url6)
$ws.Hyperlinks.Add($ws.Cells.Item(7, 6), $url7)
$ws.Hyperlinks.Add($ws.Cells.Item(8, 6), $url8)
$ws.Hyperlinks.Add($ws.Cells.Item(9, 6), $url9)
$ws.Hyperlinks.Add($ws.Cells.Item(10, 6), $url10)
$ws.Hyperlinks.Add($ws.Cells.Item(11, 6), $url11)
$ws.Hyperlinks.Add($ws.Cells.Item(12, 6), $url12)
$ws.Hyperlinks.Add($ws.Cells.Item(13, 6), $url13)

# --- Column-width tweaks: B 52->51, H 17->14 (ColumnWidth input needs a
#     -5/6 correction to land on the intended whole-character stored
#     width). ---
$ws.Columns.Item(2).ColumnWidth = 51 - 5/6
$ws.Columns.Item(8).ColumnWidth = 14 - 5/6
